# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.039.63"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "3.384.14"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "571.10"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").Value = "141.05"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").Value = "3.962.11"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").Value = "27.86"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").Value = "3.386.18"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").Value = "61.122.50"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("E18").Value = "  -1.35%  "

$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  -2.23%  "

$ws.Range("D20").Value = "8.90"
$ws.Range("E20").Value = "  -1.17%  "

$ws.Range("D21").Value = "382.24"
$ws.Range("E21").Value = "  -1.13%  "

$ws.Range("D22").Value = "75.93"
$ws.Range("E22").Value = "  +3.72%  "

$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("D26").Value = "3.518.19"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "0.189"
$ws.Range("E27").Value = "  +6.19%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").Value = "7.96"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("D34").Value = "23.28"
$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  +1.04%  "

$ws.Range("D36").Value = "166.22"
$ws.Range("E36").Value = "  -0.68%  "

$ws.Range("D37").Value = "3.419.17"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "4.98"
$ws.Range("E38").Value = "  +1.66%  "

$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("D40").Value = "0.0768"
$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("D41").Value = "26.45"
$ws.Range("E41").Value = "  -2.15%  "

$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").Value = "0.780"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("E46").Value = "  +0.67%  "

$ws.Range("D47").Value = "2.448.33"
$ws.Range("E47").Value = "  -2.65%  "

$ws.Range("D48").Value = "22.97"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("D49").Value = "6.64"
$ws.Range("E49").Value = "  -1.92%  "

$ws.Range("D50").Value = "2.13"
$ws.Range("E50").Value = "  +10.18%  "

$ws.Range("D51").Value = "0.0262"
$ws.Range("E51").Value = "  -2.15%  "
